$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new rows of issues (columns B/C/D first, in row order)
$ws.Range("B5").Value = "xóa project khi role admin"

$ws.Range("B6").Value = "kiểm tra email có thật không"
$ws.Range("C6").Value = "user"
$ws.Range("D6").Value = "cao"

$ws.Range("B7").Value = "tích hợp trello"
$ws.Range("C7").Value = "project"
$ws.Range("D7").Value = "trung bình"

# Status column: the new "đã thêm" entry first...
$ws.Range("E6").Value = "đã thêm"

# ...then the "chờ" status for the remaining rows
$ws.Range("E2").Value = "chờ"
$ws.Range("E3").Value = "chờ"
$ws.Range("E4").Value = "chờ"
$ws.Range("E5").Value = "chờ"
$ws.Range("E7").Value = "chờ"

# Match the author's final selection/cursor position
$ws.Range("E4").Select()
